$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 327, shifting existing rows (327..434) down to (328..435)
$ws.Rows("327:327").Insert()

# Populate the newly inserted row 327 with the new record
$ws.Range("A327").Value = 9
$ws.Range("B327").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C327").Value = "Metropolitana"
$ws.Range("D327").Value = 45120
$ws.Range("E327").Value = 13
$ws.Range("F327").Value = 100112043
$ws.Range("G327").Value = "Pepino ensalada"
$ws.Range("H327").Value = "Sin especificar"
$ws.Range("I327").Value = "Primera"
$ws.Range("J327").Value = 70
$ws.Range("K327").Value = 13000
$ws.Range("L327").Value = 14000
$ws.Range("M327").Value = 13500
$ws.Range("N327").Value = "`$/caja 60 unidades"
$ws.Range("O327").Value = "Región de Arica y Parinacota"
$ws.Range("P327").Value = 225
$ws.Range("Q327").Value = 60
$ws.Range("R327").Value = "Hortaliza"
